$d = $word.ActiveDocument

# Replace the body of the "Aim:" sentence, keeping the leading space and the
# trailing period in their own (unchanged) runs, matching how Word leaves
# runs split when only part of a run's text is edited.
$r = $d.Content.Duplicate
$r.Find.Execute("Choosing groups and getting started on the final project")
$r.Text = "Final project introduction, and choosing student groups"
$r.Font.Name = "Segoe UI"
$r.Font.NameFarEast = "Times New Roman"
$r.Font.NameBi = "Segoe UI"

# Re-touch the now-orphaned leading space so it keeps its own run instead of
# being re-merged with the preceding ":" run.
$sp = $d.Content.Duplicate
$sp.Find.Execute(":")
$sp.Collapse(0)
$sp.MoveEnd(1, 1)
$sp.Font.Name = "Segoe UI"
$sp.Font.NameFarEast = "Times New Roman"
$sp.Font.NameBi = "Segoe UI"
